$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename BaselineDateID -> BaselineDate (column C)
$ws.Range("C2").Value = "BaselineDate"

# Insert two new columns (ReleaseDate, ReportingDate) right after
# ProjectBusinessKey (column K), shifting the rest of the header right.
$ws.Range("L1:M1").EntireColumn.Insert()

$ws.Range("L2").Value = "ReleaseDate"
$ws.Range("M2").Value = "ReportingDate"

# Rename the shifted TargetDateID -> TargetDate (now column P)
$ws.Range("P2").Value = "TargetDate"
